# Updated account/Filters and fixed GMB issues
#
# 1) Insert a new "BingZoom" worksheet right before "Zoom", containing a
#    Day_DD / Month_MMM / Year_YYYY pair of date-parts tables (Dec 2019 /
#    Jan 2020).
# 2) Update the TPSEE sheet's filter-scenario row (Country/State/City).
# 3) Fix up sheet selections / the active tab (TPSEE becomes the active
#    sheet instead of Zoom).

$wb = $excel.ActiveWorkbook

# --- 1) Insert the new "BingZoom" sheet just before "Zoom" -----------------
$zoom = $wb.Worksheets("Zoom")
$bingZoom = $wb.Worksheets.Add($zoom)
$bingZoom.Name = "BingZoom"

$bingZoom.Range("A1").Value = "Day_DD"
$bingZoom.Range("B1").Value = "Month_MMM"
$bingZoom.Range("C1").Value = "Year_YYYY"
$bingZoom.Range("D1").Value = "Day_DD"
$bingZoom.Range("E1").Value = "Month_MMM"
$bingZoom.Range("F1").Value = "Year_YYYY"

$bingZoom.Range("A2").Value = 19
$bingZoom.Range("B2").Value = "December"
$bingZoom.Range("C2").Value = 2019
$bingZoom.Range("D2").Value = 21
$bingZoom.Range("E2").Value = "January"
$bingZoom.Range("F2").Value = 2020

$bingZoom.Range("E2").Select()

# --- 2) Update TPSEE filter values -----------------------------------------
$tpsee = $wb.Worksheets("TPSEE")
$tpsee.Range("C2").Value = "US"
$tpsee.Range("D2").Value = "Illinois"
$tpsee.Range("E2").Value = "Chicago"
$tpsee.Range("E2").Select()

# --- 3) Fix up the "Zoom" sheet selection -----------------------------------
$zoom.Range("A2:F3").Select()

# --- 4) TPSEE is now the active tab ----------------------------------------
$tpsee.Activate()
